$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: fill in G20 and H20 with 5 (style stays the same)
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 5

# Row 21: F21 changes style (no more green fill) and value to 5; G21/H21 get value 5
$ws.Range("G21").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 5

# Row 23: G23 gets value 5
$ws.Range("G23").Value = 5

# Row 27: F27 changes style and value to 5; G27/H27 get value 5
$ws.Range("G27").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 5

# Row 28: C28 and D28 change style and value to 5
$ws.Range("G28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5

# Update frozen pane top-left cell and active selection
$ws.Range("G23").Select()
$excel.ActiveWindow.ScrollRow = 7
